# Fixed formatting: adding borders and wrap text on header.
# This collapses the now-unused "No. of Sites ..." / accomplishment
# columns (X through AG) out of the sheet -- the lone still-relevant
# trailing column ("Status as of July 4, 2025", previously AH) slides
# left into the X position -- and then (re)applies border + bold
# formatting to the header row and a border to the data row beneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stale columns X:AG -- "Status as of July 4, 2025" (was AH)
# shifts left to become the new column X.
$ws.Range("X1:AG1").EntireColumn.Delete()

# Header cell (new X1) gets bold text + a thin border, same as the
# rest of row 1, but without the center/top alignment those use.
$header = $ws.Range("X1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1

# Data row gets a thin border (including the now-empty X2 dropdown
# cell), matching the rest of the data row's plain font.
$dataRow = $ws.Range("A2:X2")
$dataRow.Borders.LineStyle = 1
